# This script applies a row-permutation update to rows 2-13 of Sheet1.
# The original rows (2-13) are cyclically reshuffled: the new content of
# each row equals the old content that used to live at the row given by
# $mapping[row]. Columns A,B,C,E,F,G,H,I,J,K are constant across all rows
# so only D,L,M,N,O,P,Q,R,S,T need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that vary by row, for rows 2-13.
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$before = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping: new row r gets the old content that used to be at row $mapping[r]
$mapping = @{
    2  = 12
    3  = 13
    4  = 11
    5  = 6
    6  = 7
    7  = 4
    8  = 5
    9  = 10
    10 = 2
    11 = 3
    12 = 8
    13 = 9
}

for ($r = 2; $r -le 13; $r++) {
    $src = $mapping[$r]
    $srcVals = $before[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}
